$d = $word.ActiveDocument

# Locate the paragraph ending in "...main folder" so the new content can be
# appended right after it (and before the "_____" divider paragraph).
$rng = $d.Content
$found = $rng.Find.Execute("The default DB is db.sqlite3 in the main folder", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text for insertion point"
}

# Resolve the 1-based index (within $d.Paragraphs) of the paragraph that
# contains the matched text.
$paras = $d.Paragraphs
$anchorIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $pr = $paras.Item($i).Range
    if ($pr.Start -le $rng.Start -and $rng.End -le $pr.End) {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not resolve anchor paragraph index"
}

# Split off a brand new, empty paragraph right after the anchor paragraph.
$anchorRange = $d.Paragraphs.Item($anchorIndex).Range
$anchorRange.Collapse(0)
$anchorRange.InsertParagraphAfter()

# That empty paragraph now lives at anchorIndex + 1 - leave it empty, but
# split another paragraph after it for the "Load the data..." line.
$p1 = $d.Paragraphs.Item($anchorIndex + 1)
$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item($anchorIndex + 2)
$p2.Range.InsertBefore("Load the data with a script---")

# Split once more for the final command line.
$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item($anchorIndex + 3)
$p3.Range.InsertBefore("python3 manage.py data_to_be_loaded.py")
